$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.517.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.00%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.529.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.51%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.15%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.60%  "

# Row 7
$ws.Range("E7").Value = "  +0.31%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.506"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.42%  "

# Row 9
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.525.58"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.67%  "

# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.95%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.168"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.48%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.342"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.17%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.10%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.999.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.38%  "

# Row 15
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "70.522.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.92%  "

# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000179"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.15%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.03%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.542.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.85%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.47%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "357.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.78%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -10.27%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.06%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.77%  "

# Row 24
$ws.Range("E24").Value = "  +0.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.62%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.76%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.14%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.669.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.23%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.14%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0910"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.35%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.25%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "478.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.58%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.93%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.00%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.17%  "

# Row 36
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.44%  "

# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.115"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.17%  "

# Row 38
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.20%  "

# Row 39
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.53%  "

# Row 40
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.04%  "

# Row 41
$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.79%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.17%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.80%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.316"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.05%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.98%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.83%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.66%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.60%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.523"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.95%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.49%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.593"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.32%  "
